$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5039156666666667
$ws.Range("H2").Value = 1.511747
$ws.Range("O2").Value = 0.8416031693647025
$ws.Range("P2").Value = 0.8416031693647025
$ws.Range("Q2").Value = 0.7958743256199998
$ws.Range("R2").Value = 7.162868930579999
$ws.Range("S2").Value = 0.8416031693647025
$ws.Range("T2").Value = 0.8416031693647025

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5039156666666667
$ws.Range("H3").Value = 1.511747
$ws.Range("M3").Value = 0.2972526666666667
$ws.Range("N3").Value = 0.8917580000000001
$ws.Range("O3").Value = 0.1583968306352975
$ws.Range("P3").Value = 0.1583968306352975
$ws.Range("Q3").Value = 0.1497902756917778
$ws.Range("R3").Value = 1.348112481226
$ws.Range("S3").Value = 0.1583968306352975
$ws.Range("T3").Value = 0.1583968306352975
